$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: id=8, Nombre="Computador ", Seccion="Tecnología", Precio Compra=300000, Precio Venta=500000
$ws.Cells.Item(3, 1).Value = 8
$ws.Cells.Item(3, 2).Value = "Computador "
$ws.Cells.Item(3, 3).Value = "Tecnología"
$ws.Cells.Item(3, 4).Value = 300000
$ws.Cells.Item(3, 5).Value = 500000

# Row 4: id=9, Nombre="iPhone", Seccion="Tecnología", Precio Compra=4000000, Precio Venta=7000000
$ws.Cells.Item(4, 1).Value = 9
$ws.Cells.Item(4, 2).Value = "iPhone"
$ws.Cells.Item(4, 3).Value = "Tecnología"
$ws.Cells.Item(4, 4).Value = 4000000
$ws.Cells.Item(4, 5).Value = 7000000

# Row 5 (new): id=10, Nombre="Comedor 4 Puestos", Seccion="Hogar", Precio Compra=250000, Precio Venta=300000
$ws.Cells.Item(5, 1).Value = 10
$ws.Cells.Item(5, 2).Value = "Comedor 4 Puestos"
$ws.Cells.Item(5, 3).Value = "Hogar"
$ws.Cells.Item(5, 4).Value = 250000
$ws.Cells.Item(5, 5).Value = 300000
